$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.113.61'
$ws.Range("E2").Value = '  -0.73%  '

$ws.Range("D3").Value = '1.823.55'
$ws.Range("E3").Value = '  -0.87%  '

$ws.Range("E4").Value = '  +0.37%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '234.22'
$ws.Range("E5").Value = '  -2.04%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5968'
$ws.Range("E6").Value = '  -4.36%  '

$ws.Range("E7").Value = '  +0.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.06927'
$ws.Range("E8").Value = '  -6.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2742'
$ws.Range("E9").Value = '  -4.89%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.06'
$ws.Range("E10").Value = '  -6.65%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07604'
$ws.Range("E11").Value = '  -1.49%  '

$ws.Range("D12").Value = '1.825.30'
$ws.Range("E12").Value = '  -0.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.724'
$ws.Range("E13").Value = '  -4.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6208'
$ws.Range("E14").Value = '  -6.40%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009604'
$ws.Range("E15").Value = '  -8.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '77.70'
$ws.Range("E16").Value = '  -4.45%  '

$ws.Range("D17").Value = '28.655.58'
$ws.Range("E17").Value = '  -2.24%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.567'
$ws.Range("E18").Value = '  -10.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '216.10'
$ws.Range("E19").Value = '  -7.56%  '

$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.48'
$ws.Range("E21").Value = '  -6.13%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.810'
$ws.Range("E22").Value = '  -6.39%  '

$ws.Range("E23").Value = '  +0.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '156.51'
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '7.907'
$ws.Range("E25").Value = '  -6.35%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1278'
$ws.Range("E26").Value = '  -4.13%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.41'
$ws.Range("E27").Value = '  -4.78%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.409'
$ws.Range("E28").Value = '  -4.34%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06262'
$ws.Range("E29").Value = '  -11.52%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.436'
$ws.Range("E30").Value = '  -2.91%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.820'
$ws.Range("E31").Value = '  -4.91%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.741'
$ws.Range("E32").Value = '  -7.29%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.714'
$ws.Range("E33").Value = '  -5.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.081'
$ws.Range("E34").Value = '  -6.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6418'
$ws.Range("E35").Value = '  -8.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.539'
$ws.Range("E36").Value = '  -1.73%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.746'
$ws.Range("E37").Value = '  -1.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01745'
$ws.Range("E38").Value = '  -4.46%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.525'
$ws.Range("E39").Value = '  -3.84%  '

$ws.Range("D40").Value = '1.145.40'
$ws.Range("E40").Value = '  -6.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8816'
$ws.Range("E41").Value = '  -6.38%  '

$ws.Range("E42").Value = '  +0.39%  '

$ws.Range("D43").Value = '1.975.51'
$ws.Range("E43").Value = '  -0.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.24'
$ws.Range("E44").Value = '  -0.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '61.61'
$ws.Range("E45").Value = '  -5.11%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000112'
$ws.Range("E46").Value = '  -3.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.593'
$ws.Range("E47").Value = '  -5.30%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.432'
$ws.Range("E48").Value = '  -4.96%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05514'
$ws.Range("E49").Value = '  -2.29%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4541'
$ws.Range("E50").Value = '  -0.65%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.402'
$ws.Range("E51").Value = '  -7.63%  '
